$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This pushes the existing
# Username/Password table from A:B to B:C, matching the diff.
[void]$ws.Columns.Item(1).Insert()

# Fill the new "Tests" label column (TestNG test-case names).
$ws.Range("A1").Value = "Tests"
$ws.Range("A2").Value = "Invalid 1"
$ws.Range("A3").Value = "Invalid 2"
$ws.Range("A4").Value = "Invalid 3"
$ws.Range("A5").Value = "Valid 1"

# Rename the worksheet to match the new data purpose.
$ws.Name = "Login"

# Column widths for the new layout (A = "Tests" labels, B = Username values).
$ws.Columns.Item(1).ColumnWidth = 12.75
$ws.Columns.Item(2).ColumnWidth = 9.95

# Update the active cell / selection.
[void]$ws.Range("D3").Select()

Write-Output "done"
